$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the look of the
# existing header cells (bold font, centered/top aligned, thin box border).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

foreach ($addr in @("I1", "J1")) {
    $hdr = $ws.Range($addr)
    $hdr.Font.Bold = $true
    $hdr.HorizontalAlignment = -4108  # xlCenter
    $hdr.VerticalAlignment = -4160    # xlTop
    $hdr.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
    $hdr.Borders.Item(8).LineStyle = 1   # xlEdgeTop
    $hdr.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
    $hdr.Borders.Item(10).LineStyle = 1  # xlEdgeRight
}

# Fill in the new data columns I (I0) and J (IF) for rows 2-22
$values = @(
    @(7, 8),
    @(10, 10),
    @(8, 8),
    @(8, 9),
    @(7, 7),
    @(6, 8),
    @(6, 8),
    @(5, 7),
    @(8, 9),
    @(6, 8),
    @(7, 8),
    @(9, 9),
    @(5, 6),
    @(6, 7),
    @(7, 7),
    @(7, 8),
    @(9, 9),
    @(5, 6),
    @(6, 6),
    @(6, 6),
    @(4, 4)
)

$row = 2
foreach ($pair in $values) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
